$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 7692649
$ws.Range("I33").Value = 9091294
$ws.Range("J33").Value = 98.5
$ws.Range("K33").Value = 9091294
$ws.Range("L33").Value = 98.5
$ws.Range("M33").Value = -9091065
$ws.Range("N33").Value = -556.5
$ws.Range("H57").Value = 44234.5
$ws.Range("J57").Value = 44234.5
$ws.Range("L57").Value = 132703.5
$ws.Range("N57").Value = -133701.5
$ws.Range("H63").Value = 88979.336
$ws.Range("J63").Value = 88979.336
$ws.Range("L63").Value = 88979.336
$ws.Range("N63").Value = -90227.336
$ws.Range("H66").Value = 88979.336
$ws.Range("J66").Value = 88979.336
$ws.Range("L66").Value = 266938.008
$ws.Range("N66").Value = -273178.008
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("H69").Value = 5824
$ws.Range("J69").Value = 6780
$ws.Range("L69").Value = 20340
$ws.Range("N69").Value = -22088
$ws.Range("H70").Value = 2720.4
$ws.Range("J70").Value = 2975
$ws.Range("L70").Value = 8925
$ws.Range("N70").Value = -9465
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("H72").Value = 5824
$ws.Range("J72").Value = 6780
$ws.Range("L72").Value = 61020
$ws.Range("N72").Value = -69756
$ws.Range("H73").Value = 2720.4
$ws.Range("J73").Value = 2975
$ws.Range("L73").Value = 8925
$ws.Range("N73").Value = -10797
$ws.Range("H76").Value = 4484
$ws.Range("I76").Value = 4484
$ws.Range("K76").Value = 4484
$ws.Range("M76").Value = -4169
$ws.Range("H79").Value = 4484
$ws.Range("I79").Value = 4484
$ws.Range("K79").Value = 4484
$ws.Range("M79").Value = -3392
$ws.Range("H132").Value = 2126.303
$ws.Range("I132").Value = 1731
$ws.Range("K132").Value = 5193
$ws.Range("M132").Value = -2663
$ws.Range("H137").Value = 1779.5
$ws.Range("I137").Value = 1870.8572
$ws.Range("J137").Value = 1566.3334
$ws.Range("K137").Value = 5612.571599999999
$ws.Range("L137").Value = 4699.0002
$ws.Range("M137").Value = -3062.571599999999
$ws.Range("N137").Value = -9799.0002
$ws.Range("H141").Value = 2043.3529
$ws.Range("I141").Value = 1849.1333
$ws.Range("K141").Value = 5547.3999
$ws.Range("M141").Value = -367.3999000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 40004820
$ws.Range("I2").Value = 45459292
$ws.Range("J2").Value = 5381
$ws.Range("K2").Value = 45459292
$ws.Range("L2").Value = 5381
$ws.Range("M2").Value = -45459179
$ws.Range("N2").Value = -5607
$ws.Range("H61").Value = 2981.1462
$ws.Range("I61").Value = 2689.5405
$ws.Range("J61").Value = 5678.5
$ws.Range("K61").Value = 2689.5405
$ws.Range("L61").Value = 5678.5
$ws.Range("M61").Value = -2477.5405
$ws.Range("N61").Value = -6102.5
$ws.Range("H110").Value = 6504.64
$ws.Range("I110").Value = 6985.0527
$ws.Range("J110").Value = 4983.3335
$ws.Range("K110").Value = 6985.0527
$ws.Range("L110").Value = 4983.3335
$ws.Range("M110").Value = -4940.0527
$ws.Range("N110").Value = -9073.333500000001
$ws.Range("H116").Value = 40004820
$ws.Range("I116").Value = 45459292
$ws.Range("J116").Value = 5381
$ws.Range("K116").Value = 45459292
$ws.Range("L116").Value = 5381
$ws.Range("M116").Value = -45456998
$ws.Range("N116").Value = -9969
$ws.Range("H132").Value = 3008
$ws.Range("I132").Value = 2821.6562
$ws.Range("J132").Value = 4200.6
$ws.Range("K132").Value = 8464.9686
$ws.Range("L132").Value = 12601.8
$ws.Range("M132").Value = -5934.9686
$ws.Range("N132").Value = -17661.8
$ws.Range("H136").Value = 2981.1462
$ws.Range("I136").Value = 2689.5405
$ws.Range("J136").Value = 5678.5
$ws.Range("K136").Value = 8068.6215
$ws.Range("L136").Value = 17035.5
$ws.Range("M136").Value = -5518.6215
$ws.Range("N136").Value = -22135.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 40004820
$ws.Range("I3").Value = 45459292
$ws.Range("J3").Value = 5381
$ws.Range("K3").Value = 45459292
$ws.Range("L3").Value = 5381
$ws.Range("M3").Value = -45459178
$ws.Range("N3").Value = -5609
$ws.Range("H26").Value = 30000
$ws.Range("I26").Value = 30000
$ws.Range("K26").Value = 30000
$ws.Range("M26").Value = -29708
$ws.Range("H86").Value = 1353.7241
$ws.Range("I86").Value = 1521.7826
$ws.Range("K86").Value = 1521.7826
$ws.Range("M86").Value = -398.7826
$ws.Range("H89").Value = 1353.7241
$ws.Range("I89").Value = 1521.7826
$ws.Range("K89").Value = 7608.913
$ws.Range("M89").Value = -1992.913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23089.146
$ws.Range("I31").Value = 30943.117
$ws.Range("K31").Value = 30943.117
$ws.Range("M31").Value = -30648.117
$ws.Range("H34").Value = 23089.146
$ws.Range("I34").Value = 30943.117
$ws.Range("K34").Value = 30943.117
$ws.Range("M34").Value = -30741.117
$ws.Range("H63").Value = 60000
$ws.Range("J63").Value = 60000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61372
$ws.Range("H66").Value = 60000
$ws.Range("J66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -186864
$ws.Range("H132").Value = 3976.5134
$ws.Range("I132").Value = 3761.4333
$ws.Range("K132").Value = 11284.2999
$ws.Range("M132").Value = -8754.2999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 159.42857
$ws.Range("I10").Value = 169.33333
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 507.99999
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = -368.99999
$ws.Range("N10").Value = -578
$ws.Range("H82").Value = 28460.77
$ws.Range("H85").Value = 28460.77
$ws.Range("H106").Value = 11781.84
$ws.Range("I106").Value = 7800.5
$ws.Range("K106").Value = 23401.5
$ws.Range("M106").Value = -22455.5
$ws.Range("H131").Value = 30988.057
$ws.Range("J131").Value = 2487.7058
$ws.Range("L131").Value = 7463.117400000001
$ws.Range("N131").Value = -17543.1174

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I102").Value = 1078.3889
$ws.Range("K102").Value = 1078.3889
$ws.Range("M102").Value = 543.6111000000001
$ws.Range("I132").Value = 2726.3076
$ws.Range("J132").Value = 4115.25
$ws.Range("K132").Value = 8178.9228
$ws.Range("L132").Value = 12345.75
$ws.Range("M132").Value = -5648.9228
$ws.Range("N132").Value = -17405.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1977.0454
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 2054.75
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 2054.75
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -2644.75
$ws.Range("H27").Value = 1977.0454
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 2054.75
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 2054.75
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -2268.75
$ws.Range("H29").Value = 3849.3333
$ws.Range("J29").Value = 3849.3333
$ws.Range("L29").Value = 3849.3333
$ws.Range("N29").Value = -4439.3333
$ws.Range("H68").Value = 3886.182
$ws.Range("I68").Value = 2999.8333
$ws.Range("K68").Value = 2999.8333
$ws.Range("M68").Value = -2250.8333
$ws.Range("H71").Value = 3886.182
$ws.Range("I71").Value = 2999.8333
$ws.Range("K71").Value = 14999.1665
$ws.Range("M71").Value = -11255.1665
$ws.Range("H82").Value = 2099.4211
$ws.Range("I82").Value = 1971.3572
$ws.Range("K82").Value = 1971.3572
$ws.Range("M82").Value = -1610.3572
$ws.Range("H85").Value = 2099.4211
$ws.Range("I85").Value = 1971.3572
$ws.Range("K85").Value = 1971.3572
$ws.Range("M85").Value = -723.3571999999999
$ws.Range("H136").Value = 8000.5
$ws.Range("I136").Value = 7500.6665
$ws.Range("K136").Value = 22501.9995
$ws.Range("M136").Value = -19951.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2576.625
$ws.Range("I132").Value = 2558.2173
$ws.Range("K132").Value = 7674.651899999999
$ws.Range("M132").Value = -5144.651899999999
$ws.Range("H136").Value = 3454.5454
$ws.Range("J136").Value = 3428.5715
$ws.Range("L136").Value = 10285.7145
$ws.Range("N136").Value = -15385.7145
